# Rename the three worksheets to their Hebrew names and make the third
# sheet ("מכפלה", formerly "Sheet3") the active/selected sheet.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet1").Name = "כללי"
$wb.Worksheets.Item("Sheet2").Name = "מכרז חפץ אחד"
$wb.Worksheets.Item("Sheet3").Name = "מכפלה"

$wb.Worksheets.Item("מכפלה").Activate()
